# TimeSheet.xlsx update — "Updated Timesheet till 9/2/2012"
#
# - Jan: days 1/25-1/28 flip from "NA / OFF" (non-working) to worked days
#        logging P01 activity, 1 hour each.
# - Feb: days 2/1-2/9 get their Activity Code / Hours logged for the first time.
# - View state: Apr becomes the active tab/sheet; each sheet's remembered
#   selection moves on to reflect where the user was last working.
# - Column J (the "Activity" description column) is widened slightly on every
#   sheet to better fit its text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Jan sheet: rows 31-34 (2013-01-25 .. 2013-01-28)
#   C: "NA" -> "P01"   (still text)
#   D: "OFF" (text) -> 1 (numeric hours)
# ---------------------------------------------------------------------------
$jan = $wb.Worksheets.Item("Jan")

$jan.Range("C31").Value = "P01"
$jan.Range("D31").Value = 1

$jan.Range("C32").Value = "P01"
$jan.Range("D32").Value = 1

$jan.Range("C33").Value = "P01"
$jan.Range("D33").Value = 1

$jan.Range("C34").Value = "P01"
$jan.Range("D34").Value = 1

# ---------------------------------------------------------------------------
# Feb sheet: rows 7-15 (2013-02-01 .. 2013-02-09) get new Activity/Hours data
# ---------------------------------------------------------------------------
$feb = $wb.Worksheets.Item("Feb")

$feb.Range("C7").Value = "P01"
$feb.Range("D7").Value = 1

$feb.Range("C8").Value = "P02"
$feb.Range("D8").Value = 1

$feb.Range("C9").Value = "P02"
$feb.Range("D9").Value = 1

$feb.Range("C10").Value = "P03"
$feb.Range("D10").Value = 0.5

$feb.Range("C11").Value = "NA"
$feb.Range("D11").Value = 2

$feb.Range("C12").Value = "P03"
$feb.Range("D12").Value = "OFF"

$feb.Range("C13").Value = "NA"
$feb.Range("D13").Value = "OFF"

$feb.Range("C14").Value = "NA"
$feb.Range("D14").Value = "OFF"

$feb.Range("C15").Value = "P03"
$feb.Range("D15").Value = 1

# ---------------------------------------------------------------------------
# Column J width tweak (col 10) on every sheet — each sheet gets its own
# slightly different width so the longest "Activity" description fits.
# ---------------------------------------------------------------------------
$jan.Columns.Item(10).ColumnWidth = 61.666666666666664
$feb.Columns.Item(10).ColumnWidth = 70.5

$mar = $wb.Worksheets.Item("Mar")
$mar.Columns.Item(10).ColumnWidth = 68.33333333333333

$apr = $wb.Worksheets.Item("Apr")
$apr.Columns.Item(10).ColumnWidth = 73

# ---------------------------------------------------------------------------
# View state: walk the sheets in tab order, moving the remembered selection
# on each one, finishing on Apr so it ends up the active tab/sheet.
# ---------------------------------------------------------------------------
$jan.Activate()
$jan.Range("C36").Select()

$feb.Activate()
$feb.Range("D10").Select()

$mar.Activate()
$mar.Range("B24").Select()

$apr.Activate()
$apr.Range("E15").Select()
